$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 1701.1538
$ws.Range("I39").Value = 146.33333
$ws.Range("K39").Value = 438.99999
$ws.Range("M39").Value = -142.99999
$ws.Range("H54").Value = 38000
$ws.Range("I54").Value = 38000
$ws.Range("K54").Value = 38000
$ws.Range("M54").Value = -37514
$ws.Range("H58").Value = 1015
$ws.Range("I58").Value = 1015
$ws.Range("K58").Value = 3045
$ws.Range("M58").Value = -2895
$ws.Range("H80").Value = 1030.3077
$ws.Range("I80").Value = 320.33334
$ws.Range("J80").Value = 1638.8572
$ws.Range("K80").Value = 961.0000200000001
$ws.Range("L80").Value = 4916.571599999999
$ws.Range("M80").Value = 36.99997999999994
$ws.Range("N80").Value = -6912.571599999999
$ws.Range("H83").Value = 1030.3077
$ws.Range("I83").Value = 320.33334
$ws.Range("J83").Value = 1638.8572
$ws.Range("K83").Value = 2883.00006
$ws.Range("L83").Value = 14749.7148
$ws.Range("M83").Value = 2108.99994
$ws.Range("N83").Value = -24733.7148
$ws.Range("H92").Value = 1484.96
$ws.Range("I92").Value = 357.29413
$ws.Range("J92").Value = 3881.25
$ws.Range("K92").Value = 357.29413
$ws.Range("L92").Value = 3881.25
$ws.Range("M92").Value = 890.70587
$ws.Range("N92").Value = -6377.25
$ws.Range("H103").Value = 1406.5834
$ws.Range("I103").Value = 1955.1428
$ws.Range("K103").Value = 5865.428400000001
$ws.Range("M103").Value = -5279.428400000001
$ws.Range("H106").Value = 5700669.5
$ws.Range("I106").Value = 6510951
$ws.Range("J106").Value = 28699.5
$ws.Range("K106").Value = 6510951
$ws.Range("L106").Value = 28699.5
$ws.Range("M106").Value = -6510320
$ws.Range("N106").Value = -29961.5
$ws.Range("H111").Value = 1944.2273
$ws.Range("I111").Value = 1527.1428
$ws.Range("J111").Value = 2674.125
$ws.Range("K111").Value = 4581.428400000001
$ws.Range("L111").Value = 8022.375
$ws.Range("M111").Value = -1514.428400000001
$ws.Range("N111").Value = -14156.375
$ws.Range("H116").Value = 10828.733
$ws.Range("J116").Value = 9082.223
$ws.Range("L116").Value = 9082.223
$ws.Range("N116").Value = -15966.223
$ws.Range("H129").Value = 1840.8096
$ws.Range("I129").Value = 527.5833
$ws.Range("J129").Value = 3591.7778
$ws.Range("K129").Value = 1582.7499
$ws.Range("L129").Value = 10775.3334
$ws.Range("M129").Value = 3417.2501
$ws.Range("N129").Value = -20775.3334
$ws.Range("H132").Value = 3248.7646
$ws.Range("I132").Value = 3045.3
$ws.Range("K132").Value = 9135.900000000001
$ws.Range("M132").Value = -6605.900000000001
$ws.Range("H138").Value = 3290.5417
$ws.Range("I138").Value = 3429.889
$ws.Range("J138").Value = 3206.9333
$ws.Range("K138").Value = 10289.667
$ws.Range("L138").Value = 9620.7999
$ws.Range("M138").Value = -5149.667000000001
$ws.Range("N138").Value = -19900.7999
$ws.Range("H141").Value = 5821
$ws.Range("I141").Value = 3501.7896
$ws.Range("J141").Value = 16837.25
$ws.Range("K141").Value = 10505.3688
$ws.Range("L141").Value = 50511.75
$ws.Range("M141").Value = -5325.3688
$ws.Range("N141").Value = -60871.75

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H22").Value = 10000
$ws.Range("J22").Value = 10000
$ws.Range("L22").Value = 10000
$ws.Range("N22").Value = -10598
$ws.Range("H30").Value = 6960
$ws.Range("J30").Value = 7000
$ws.Range("L30").Value = 7000
$ws.Range("N30").Value = -7300
$ws.Range("H32").Value = 3212.82
$ws.Range("I32").Value = 2744.0476
$ws.Range("K32").Value = 2744.0476
$ws.Range("M32").Value = -2457.0476
$ws.Range("H45").Value = 2062.375
$ws.Range("I45").Value = 1105.3846
$ws.Range("K45").Value = 1105.3846
$ws.Range("M45").Value = -728.3846000000001
$ws.Range("H74").Value = 1778.0588
$ws.Range("I74").Value = 1041.0526
$ws.Range("J74").Value = 2711.6
$ws.Range("K74").Value = 1041.0526
$ws.Range("L74").Value = 2711.6
$ws.Range("M74").Value = -167.0526
$ws.Range("N74").Value = -4459.6
$ws.Range("H77").Value = 1778.0588
$ws.Range("I77").Value = 1041.0526
$ws.Range("J77").Value = 2711.6
$ws.Range("K77").Value = 5205.263
$ws.Range("L77").Value = 13558
$ws.Range("M77").Value = -837.2629999999999
$ws.Range("N77").Value = -22294
$ws.Range("H96").Value = 31172
$ws.Range("J96").Value = 31172
$ws.Range("L96").Value = 31172
$ws.Range("N96").Value = -36664
$ws.Range("H102").Value = 2679
$ws.Range("I102").Value = 2347.4285
$ws.Range("K102").Value = 2347.4285
$ws.Range("M102").Value = -725.4285

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 36466.168
$ws.Range("I26").Value = 31759.4
$ws.Range("J26").Value = 60000
$ws.Range("K26").Value = 31759.4
$ws.Range("L26").Value = 60000
$ws.Range("M26").Value = -31467.4
$ws.Range("N26").Value = -60584
$ws.Range("H86").Value = 9823.333000000001
$ws.Range("J86").Value = 14433.4
$ws.Range("L86").Value = 14433.4
$ws.Range("N86").Value = -16679.4
$ws.Range("H89").Value = 9823.333000000001
$ws.Range("J89").Value = 14433.4
$ws.Range("L89").Value = 72167
$ws.Range("N89").Value = -83399
$ws.Range("H134").Value = 8449.290000000001
$ws.Range("I134").Value = 7147.3477
$ws.Range("J134").Value = 12192.375
$ws.Range("K134").Value = 21442.0431
$ws.Range("L134").Value = 36577.125
$ws.Range("M134").Value = -18907.0431
$ws.Range("N134").Value = -41647.125

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 16815.875
$ws.Range("I62").Value = 4247.6665
$ws.Range("K62").Value = 4247.6665
$ws.Range("M62").Value = -3623.6665
$ws.Range("H65").Value = 16815.875
$ws.Range("I65").Value = 4247.6665
$ws.Range("K65").Value = 21238.3325
$ws.Range("M65").Value = -18118.3325
$ws.Range("H94").Value = 3303.6428
$ws.Range("I94").Value = 4988.5713
$ws.Range("J94").Value = 1618.7142
$ws.Range("K94").Value = 4988.5713
$ws.Range("L94").Value = 1618.7142
$ws.Range("M94").Value = -4537.5713
$ws.Range("N94").Value = -2520.7142
$ws.Range("H105").Value = 1283.7142
$ws.Range("I105").Value = 1161.125
$ws.Range("J105").Value = 1676
$ws.Range("K105").Value = 1161.125
$ws.Range("L105").Value = 1676
$ws.Range("M105").Value = 585.875
$ws.Range("N105").Value = -5170
$ws.Range("H112").Value = 108138.2
$ws.Range("J112").Value = 108138.2
$ws.Range("L112").Value = 108138.2
$ws.Range("N112").Value = -111092.2
$ws.Range("H132").Value = 5069.857
$ws.Range("I132").Value = 10200
$ws.Range("J132").Value = 3017.8
$ws.Range("K132").Value = 30600
$ws.Range("L132").Value = 9053.400000000001
$ws.Range("M132").Value = -28070
$ws.Range("N132").Value = -14113.4

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 3577.375
$ws.Range("J114").Value = 4488.1665
$ws.Range("L114").Value = 13464.4995
$ws.Range("N114").Value = -19972.4995
$ws.Range("H129").Value = 27781204
$ws.Range("J129").Value = 33337302
$ws.Range("L129").Value = 100011906
$ws.Range("N129").Value = -100021906
$ws.Range("H131").Value = 5232479
$ws.Range("I131").Value = 18521102
$ws.Range("J131").Value = 4223216.5
$ws.Range("K131").Value = 55563306
$ws.Range("L131").Value = 12669649.5
$ws.Range("M131").Value = -55558266
$ws.Range("N131").Value = -12679729.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 460.42856
$ws.Range("I2").Value = 422.36365
$ws.Range("K2").Value = 422.36365
$ws.Range("M2").Value = -309.36365
$ws.Range("H70").Value = 7200.609
$ws.Range("J70").Value = 9249.25
$ws.Range("L70").Value = 9249.25
$ws.Range("N70").Value = -9789.25
$ws.Range("H73").Value = 7200.609
$ws.Range("J73").Value = 9249.25
$ws.Range("L73").Value = 9249.25
$ws.Range("N73").Value = -11121.25
$ws.Range("H122").Value = 2477.1428
$ws.Range("J122").Value = 2617.75
$ws.Range("L122").Value = 7853.25
$ws.Range("N122").Value = -12753.25
$ws.Range("H126").Value = 6848.3
$ws.Range("I126").Value = 6199.6665
$ws.Range("K126").Value = 18598.9995
$ws.Range("M126").Value = -16128.9995

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 10849.667
$ws.Range("I16").Value = 10849.667
$ws.Range("K16").Value = 10849.667
$ws.Range("M16").Value = -10679.667
$ws.Range("H46").Value = 1495.55
$ws.Range("I46").Value = 940.9
$ws.Range("K46").Value = 940.9
$ws.Range("M46").Value = -752.9
$ws.Range("H136").Value = 2236.8333
$ws.Range("I136").Value = 1117.0769
$ws.Range("J136").Value = 3560.182
$ws.Range("K136").Value = 3351.2307
$ws.Range("L136").Value = 10680.546
$ws.Range("M136").Value = -801.2307000000001
$ws.Range("N136").Value = -15780.546

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 59761.285
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 59761.285
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 59761.285
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -60743.285
$ws.Range("H100").Value = 972.4666999999999
$ws.Range("I100").Value = 331.81818
$ws.Range("J100").Value = 2734.25
$ws.Range("K100").Value = 663.63636
$ws.Range("L100").Value = 5468.5
$ws.Range("M100").Value = -122.63636
$ws.Range("N100").Value = -6550.5
$ws.Range("H107").Value = 5010.174
$ws.Range("J107").Value = 1000
$ws.Range("L107").Value = 3000
$ws.Range("N107").Value = -6840
$ws.Range("H139").Value = 70000
$ws.Range("J139").Value = 70000
$ws.Range("L139").Value = 70000
$ws.Range("N139").Value = -80280
